$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "59.846.40"
Set-TextValue $ws.Range("E2") "  -1.96%  "

Set-TextValue $ws.Range("D3") "2.304.37"
Set-TextValue $ws.Range("E3") "  -2.79%  "

Set-TextValue $ws.Range("E4") "  +0.09%  "

Set-TextValue $ws.Range("D5") "540.80"
Set-TextValue $ws.Range("E5") "  -1.27%  "

Set-TextValue $ws.Range("D6") "128.61"
Set-TextValue $ws.Range("E6") "  -3.13%  "

Set-TextValue $ws.Range("E7") "  +0.06%  "

Set-TextValue $ws.Range("E8") "  -3.34%  "

Set-TextValue $ws.Range("D9") "2.302.75"
Set-TextValue $ws.Range("E9") "  -2.77%  "

Set-TextValue $ws.Range("E10") "  -0.95%  "

Set-TextValue $ws.Range("D11") "5.52"
Set-TextValue $ws.Range("E11") "  +0.25%  "

Set-TextValue $ws.Range("E12") "  -0.55%  "

Set-TextValue $ws.Range("E13") "  -2.03%  "

Set-TextValue $ws.Range("B14") "Avalanche"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D14") "23.11"
Set-TextValue $ws.Range("E14") "  -4.54%  "

Set-TextValue $ws.Range("D15") "2.714.55"
Set-TextValue $ws.Range("E15") "  -2.71%  "

Set-TextValue $ws.Range("B16") "WrappedBTC"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D16") "59.767.65"
Set-TextValue $ws.Range("E16") "  -1.86%  "

Set-TextValue $ws.Range("E17") "  -2.02%  "

Set-TextValue $ws.Range("D18") "2.291.50"
Set-TextValue $ws.Range("E18") "  -2.66%  "

Set-TextValue $ws.Range("D19") "10.40"
Set-TextValue $ws.Range("E19") "  -3.38%  "

Set-TextValue $ws.Range("E20") "  -4.36%  "

Set-TextValue $ws.Range("D21") "309.27"
Set-TextValue $ws.Range("E21") "  -2.69%  "

Set-TextValue $ws.Range("E22") "  -8.64%  "

Set-TextValue $ws.Range("E23") "  +0.04%  "

Set-TextValue $ws.Range("D24") "63.18"
Set-TextValue $ws.Range("E24") "  -0.60%  "

Set-TextValue $ws.Range("E25") "  -2.57%  "

Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.18%  "

Set-TextValue $ws.Range("D27") "7.72"
Set-TextValue $ws.Range("E27") "  -4.67%  "

Set-TextValue $ws.Range("D28") "1.35"
Set-TextValue $ws.Range("E28") "  -0.64%  "

Set-TextValue $ws.Range("E29") "  +0.38%  "

Set-TextValue $ws.Range("D30") "1.18"
Set-TextValue $ws.Range("E30") "  +2.72%  "

Set-TextValue $ws.Range("E31") "  -2.64%  "

Set-TextValue $ws.Range("E32") "  -4.17%  "

Set-TextValue $ws.Range("D33") "5.79"
Set-TextValue $ws.Range("E33") "  -2.13%  "

Set-TextValue $ws.Range("E34") "  -1.94%  "

Set-TextValue $ws.Range("B35") "ImmutableX"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D35") "1.34"
Set-TextValue $ws.Range("E35") "  -4.93%  "

Set-TextValue $ws.Range("B36") "USDe"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D36") "0.999"
Set-TextValue $ws.Range("E36") "  -0.01%  "

Set-TextValue $ws.Range("D37") "17.67"
Set-TextValue $ws.Range("E37") "  -2.52%  "

Set-TextValue $ws.Range("E38") "  +0.22%  "

Set-TextValue $ws.Range("D39") "4.02"
Set-TextValue $ws.Range("E39") "  -5.24%  "

Set-TextValue $ws.Range("D40") "310.81"
Set-TextValue $ws.Range("E40") "  -3.22%  "

Set-TextValue $ws.Range("D41") "37.93"
Set-TextValue $ws.Range("E41") "  -1.16%  "

Set-TextValue $ws.Range("E42") "  -4.48%  "

Set-TextValue $ws.Range("D43") "135.66"
Set-TextValue $ws.Range("E43") "  -5.94%  "

Set-TextValue $ws.Range("E44") "  -2.59%  "

Set-TextValue $ws.Range("D45") "0.0934"
Set-TextValue $ws.Range("E45") "  -2.14%  "

Set-TextValue $ws.Range("E46") "  -0.16%  "

Set-TextValue $ws.Range("B47") "InjectiveProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "18.46"
Set-TextValue $ws.Range("E47") "  -5.83%  "

Set-TextValue $ws.Range("B48") "Hedera"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D48") "0.0488"
Set-TextValue $ws.Range("E48") "  -2.94%  "

Set-TextValue $ws.Range("D49") "0.0₆0224"
Set-TextValue $ws.Range("E49") "  +9.94%  "

Set-TextValue $ws.Range("E50") "  -1.27%  "

Set-TextValue $ws.Range("D51") "10.97"
Set-TextValue $ws.Range("E51") "  -0.60%  "
